$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 699, shifting existing rows 699:790 down to 700:791.
$ws.Rows.Item(699).Insert()

# Populate the newly inserted row 699 with the new data record.
$ws.Range("A699").Value = 6
$ws.Range("B699").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C699").Value = "Metropolitana"
$ws.Range("D699").Value = 45131
$ws.Range("D699").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E699").Value = 13
$ws.Range("F699").Value = 100112012
$ws.Range("G699").Value = "Espinaca"
$ws.Range("H699").Value = "Sin especificar"
$ws.Range("I699").Value = "Primera"
$ws.Range("J699").Value = 350
$ws.Range("K699").Value = 7500
$ws.Range("L699").Value = 8000
$ws.Range("M699").Value = 7714
$ws.Range("N699").Value = "`$/cuna 10 kilos"
$ws.Range("O699").Value = "Región Metropolitana"
$ws.Range("P699").Value = 771
$ws.Range("Q699").Value = 10
$ws.Range("R699").Value = "Hortaliza"
